$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Resize the four columns of the "Product Backlog" table to match the
# tblGrid definition (2553 / 2199 / 2103 / 2207 dxa -> points = dxa/20).
$widths = @(127.65, 109.95, 105.15, 110.35)
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    for ($c = 1; $c -le 4; $c++) {
        $t.Cell($r, $c).Width = $widths[$c - 1]
    }
}

# Remove the trailing empty "L'utilisateur veut" row (last row of the table).
$t.Rows.Item($t.Rows.Count).Delete()

# Fix up the title: merge the "Pro" + "duct " runs into a single
# "Product " run, keeping "Backlog" untouched.
$findRange = $d.Content
$findRange.Find.Execute("Product Backlog") | Out-Null
$proStart = $findRange.Start
$proRange = $d.Range($proStart, $proStart + 3)
$proRange.Delete()

$findRange2 = $d.Content
$findRange2.Find.Execute("duct Backlog") | Out-Null
$ductStart = $findRange2.Start
$ductRange = $d.Range($ductStart, $ductStart + 5)
$ductRange.Text = "Product "
